$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextCell 'D2' '61.693.65'
Set-TextCell 'E2' '  +5.42%  '
Set-TextCell 'D3' '3.056.56'
Set-TextCell 'E3' '  +3.04%  '
Set-TextCell 'E4' '  +0.24%  '
Set-TextCell 'D5' '574.77'
Set-TextCell 'E5' '  +2.95%  '
Set-TextCell 'D6' '142.20'
Set-TextCell 'E6' '  +4.58%  '
Set-TextCell 'E7' '  +0.07%  '
Set-TextCell 'D8' '3.043.62'
Set-TextCell 'E8' '  +2.89%  '
Set-TextCell 'D9' '0.523'
Set-TextCell 'E9' '  +1.43%  '
Set-TextCell 'D10' '0.138'
Set-TextCell 'E10' '  +6.24%  '
Set-TextCell 'D11' '5.46'
Set-TextCell 'E11' '  +13.50%  '
Set-TextCell 'D12' '0.464'
Set-TextCell 'E12' '  +2.28%  '
Set-TextCell 'D13' '0.0000237'
Set-TextCell 'E13' '  +5.45%  '
Set-TextCell 'D14' '34.65'
Set-TextCell 'E15' '  -0.01%  '
Set-TextCell 'D16' '3.562.04'
Set-TextCell 'E16' '  +3.11%  '
Set-TextCell 'D17' '7.18'
Set-TextCell 'E17' '  +3.73%  '
Set-TextCell 'D18' '3.054.42'
Set-TextCell 'E18' '  +3.09%  '
Set-TextCell 'D19' '61.678.27'
Set-TextCell 'D20' '447.72'
Set-TextCell 'E20' '  +6.79%  '
Set-TextCell 'E21' '  +3.58%  '
Set-TextCell 'D22' '0.728'
Set-TextCell 'E22' '  +2.89%  '
Set-TextCell 'D23' '7.25'
Set-TextCell 'E23' '  +2.76%  '
Set-TextCell 'D24' '13.58'
Set-TextCell 'E24' '  +2.41%  '
Set-TextCell 'D25' '81.69'
Set-TextCell 'E25' '  +2.28%  '
Set-TextCell 'E26' '  +0.03%  '
Set-TextCell 'D27' '2.23'
Set-TextCell 'E27' '  +6.43%  '
Set-TextCell 'D28' '1.00'
Set-TextCell 'E28' '  +0.37%  '
Set-TextCell 'D29' '2.63'
Set-TextCell 'E29' '  +4.89%  '
Set-TextCell 'D30' '7.98'
Set-TextCell 'E30' '  +4.02%  '
Set-TextCell 'D31' '6.50'
Set-TextCell 'E31' '  +8.03%  '
Set-TextCell 'D32' '26.40'
Set-TextCell 'E32' '  +3.70%  '
Set-TextCell 'D33' '0.106'
Set-TextCell 'E33' '  +7.45%  '
Set-TextCell 'D34' '0.0₃0807'
Set-TextCell 'E34' '  +8.33%  '
Set-TextCell 'E35' '  +3.06%  '
Set-TextCell 'D36' '6.04'
Set-TextCell 'E36' '  +6.04%  '
Set-TextCell 'D37' '2.17'
Set-TextCell 'E37' '  +5.87%  '
Set-TextCell 'D38' '50.07'
Set-TextCell 'E38' '  +3.08%  '
Set-TextCell 'D39' '2.94'
Set-TextCell 'E39' '  +7.75%  '
Set-TextCell 'E40' '  +2.30%  '
Set-TextCell 'D41' '413.72'
Set-TextCell 'E41' '  +4.67%  '
Set-TextCell 'E42' '  +5.88%  '
Set-TextCell 'D43' '2.767.69'
Set-TextCell 'E43' '  +1.17%  '
Set-TextCell 'D44' '0.108'
Set-TextCell 'E44' '  +1.06%  '
Set-TextCell 'D45' '0.262'
Set-TextCell 'E45' '  +8.83%  '
Set-TextCell 'D46' '36.66'
Set-TextCell 'E46' '  +15.59%  '
Set-TextCell 'B47' 'USDe'
Set-TextCell 'C47' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextCell 'D47' '0.999'
Set-TextCell 'E47' '  -0.04%  '
Set-TextCell 'B48' 'Fetch.AI'
Set-TextCell 'C48' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell 'D48' '2.08'
Set-TextCell 'E48' '  +4.44%  '
Set-TextCell 'D49' '122.82'
Set-TextCell 'E49' '  -1.48%  '
Set-TextCell 'E50' '  +1.85%  '
Set-TextCell 'D51' '23.97'
Set-TextCell 'E51' '  +3.92%  '
